$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -4122 == xlPasteFormats
$xlPasteFormats = -4122

# Add new "ERROR" column header (copy header style/formatting from D1)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial($xlPasteFormats)

$ws.Range("E1").Value = "ERROR"
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 0.05936598391918613
$ws.Range("E2").Value = 0.003035493077078442
$ws.Range("E3").Value = 0.0005326809152579089
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 0.04806465453323013
$ws.Range("E4").Value = 0.002708960911056775
$ws.Range("E5").Value = 0.001911775804863956
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 0.01339031299261635
$ws.Range("E6").Value = 0.001508172562733018
$ws.Range("E7").Value = 0.001176774733969529
$ws.Range("E8").Value = 0.00186518307382544
$ws.Range("E9").Value = 0.002112894099568819
$ws.Range("E10").Value = 0.002449264765007386
$ws.Range("E11").Value = 0.00261466787634065
$ws.Range("E12").Value = 0.003501666092421833
$ws.Range("E13").Value = 0.00439895220724949
$ws.Range("E14").Value = 0.0006447551638369098
$ws.Range("E15").Value = 0.001541576703823349
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0.01572741592980416
$ws.Range("E16").Value = 0.001522266335368331
$ws.Range("E17").Value = 0.004903073700746497
$ws.Range("E18").Value = 0.0007796941141099332
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.008673820530755825
$ws.Range("E19").Value = 0.001461840480067992
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0.01461474403233252
$ws.Range("E20").Value = 0.001123036702029975
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 0.005752292219242671
$ws.Range("E21").Value = 0.0008842692873798057
$ws.Range("E22").Value = 0.002141996245434148
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0.06435848202244433
$ws.Range("E23").Value = 0.001741532107087278
$ws.Range("E24").Value = 0.003099782478848739
$ws.Range("E25").Value = 0.001743835996372207
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 0.03374938964761337
$ws.Range("E26").Value = 0.001654950577338617
$ws.Range("E27").Value = 0.001304340853869849
$ws.Range("E28").Value = 0.002070416134208245
$ws.Range("E29").Value = 0.003851816244631217
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 0.07163892755319588
$ws.Range("E30").Value = 0.00383757027733713
$ws.Range("E31").Value = 0.001886740493001429
$ws.Range("E32").Value = 0.004455258641850102
$ws.Range("E33").Value = 0.008940562847279374
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 0.02677773877005402
$ws.Range("E34").Value = 0.006987578924655492
$ws.Range("E35").Value = 0.003410544974445649
$ws.Range("E36").Value = 0.002339332834741048
$ws.Range("E37").Value = 0.001012704540073741
$ws.Range("E38").Value = 0.001801517900222171
$ws.Range("E39").Value = 0.001475055964619051
$ws.Range("E40").Value = 0.001466410306524519
$ws.Range("E41").Value = 0.002542866581186163
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 0.006157213440128234
$ws.Range("E42").Value = 0.002397460528171445
$ws.Range("E43").Value = 0.001018322376843115
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 0.005819961319551851
$ws.Range("E44").Value = 0.0009469795172086425
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 0.0269225597187095
$ws.Range("E45").Value = 0.001863736576094903
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 0.01410576666286996
$ws.Range("E46").Value = 0.001100481841451729
$ws.Range("E47").Value = 0.001227981137154653
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 0.005750758309528895
$ws.Range("E48").Value = 0
$ws.Range("B49").Value = 2491.878218119664
$ws.Range("C49").Value = 2
$ws.Range("D49").Value = [double]"5.516487026791419e-15"
$ws.Range("E49").Value = 0
$ws.Range("B50").Value = 2511.967542799203
$ws.Range("C50").Value = 3
$ws.Range("D50").Value = 0.0100360798260475
$ws.Range("E50").Value = 0.001343369766402542
$ws.Range("B51").Value = 2521.402026059814
$ws.Range("C51").Value = 3
$ws.Range("D51").Value = 0.04867836555897274
$ws.Range("E51").Value = 0.002802644367327954
$ws.Range("B52").Value = 2532.743580491942
$ws.Range("D52").Value = 0.01096972208564659
$ws.Range("E52").Value = 0.001281034991850073
$ws.Range("B53").Value = 2542.951706183542
$ws.Range("C53").Value = 3
$ws.Range("D53").Value = 0.02919350932781831
$ws.Range("E53").Value = 0.001775660183078727
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 2552.6607271736
$ws.Range("C54").Value = 3
$ws.Range("D54").Value = 0.02048890572107748
$ws.Range("E54").Value = 0.001545096733683355

# Apply the bold/bordered "index" column style to the newly added row 54, column A
$ws.Range("A2").Copy()
$ws.Range("A54").PasteSpecial($xlPasteFormats)
